$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# --- Column J ("asana task") updates: new task links for previously-blank cells ---
$ws.Range("J14").Value = "https://app.asana.com/0/0/1202626450708380/f"
$ws.Range("J15").Value = "https://app.asana.com/0/0/1202662876083700/f"
$ws.Range("J16").Value = "https://app.asana.com/0/0/1202677696411529/f"
$ws.Range("J18").Value = "https://app.asana.com/0/0/1202503525669953/f"

# --- Column I ("network number") updates ---
$ws.Range("I18").Value = "NGF_Networks_01"
$ws.Range("I14").Value = "NGF_Networks_02"
$ws.Range("I3").Value  = "TM1_2015_Base_Network"
$ws.Range("I4").Value  = "TM1_2015_Base_Network"
$ws.Range("I5").Value  = "TM1_2015_Base_Network"
$ws.Range("I6").Value  = "TM1_2015_Base_Network"
$ws.Range("I7").Value  = "BlueprintNetworks_64"
$ws.Range("I8").Value  = "existing_conditions\net_2035_NextGenFwy"
$ws.Range("I9").Value  = "NGF_Networks_Blueprint_01"
$ws.Range("I10").Value = "NGF_Networks_Blueprint_01"
$ws.Range("I11").Value = "NGF_Networks_Blueprint_01"
$ws.Range("I12").Value = "NGF_Networks_Blueprint_01"
$ws.Range("I13").Value = "NGF_Networks_Blueprint_01"
$ws.Range("I15").Value = "NGF_Networks_Blueprint_01"
$ws.Range("I16").Value = "NGF_Networks_Blueprint_01"

# --- New hyperlink on J8 (text unchanged, target matches the URL already shown) ---
$existingJ8 = $ws.Range("J8").Value()
$ws.Hyperlinks.Add($ws.Range("J8"), $existingJ8)

# --- Selection / active cell, matching the author's final cursor position ---
$ws.Range("H34").Select()
